$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weather Data")

$ws.Range("B7").Value = "12.55 °C (Feels like 10.88 °C)"
$ws.Range("B8").Value = "12.32 °C to 13.91 °C"
$ws.Range("B9").Value = "1030 hPa"

$humidity = $ws.Range("B10")
$humidity.NumberFormat = "@"
$humidity.Value = "39%"

$ws.Range("B11").Value = "1.54 m/s at 140°"
